$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.65
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 6.25
$ws.Range("J3").Value = 2.3
$ws.Range("K3").Value = 2.05
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.3
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62
$ws.Range("AA3").Value = 15
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 7
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 21
$ws.Range("AN3").Value = 3.4
$ws.Range("AO3").Value = 9
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 29
$ws.Range("AR3").Value = 51
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 9.5

# Row 4
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 2.15
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.57
$ws.Range("P4").Value = 2.25
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
$ws.Range("W4").Value = 8
$ws.Range("AC4").Value = 5.5
$ws.Range("AK4").Value = 23

# Row 5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 3.35
$ws.Range("J7").Value = 2.95
$ws.Range("K7").Value = 2.05
$ws.Range("N7").Value = 7.6
$ws.Range("O7").Value = 1.27
$ws.Range("P7").Value = 3.45
$ws.Range("Q7").Value = 1.82
$ws.Range("R7").Value = 1.93
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.6
$ws.Range("X7").Value = 11.75
$ws.Range("Z7").Value = 23
$ws.Range("AC7").Value = 7.6
$ws.Range("AD7").Value = 6.5
$ws.Range("AG7").Value = 10
$ws.Range("AL7").Value = 29
$ws.Range("AN7").Value = 4.25
$ws.Range("AP7").Value = 22
$ws.Range("AQ7").Value = 55
$ws.Range("AT7").Value = 2.6
$ws.Range("AW7").Value = 4.8
$ws.Range("AY7").Value = 24
$ws.Range("AZ7").Value = 80
$ws.Range("BA7").Value = 120
$ws.Range("BB7").Value = 350

# Row 8
$ws.Range("G8").Value = 2.02
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 3.3
$ws.Range("J8").Value = 2.57
$ws.Range("K8").Value = 2.18
$ws.Range("L8").Value = 3.8
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 8
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 3.65
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 2.02
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 2.9
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.15
$ws.Range("W8").Value = 8.75
$ws.Range("X8").Value = 10.75
$ws.Range("AA8").Value = 15
$ws.Range("AB8").Value = 22
$ws.Range("AC8").Value = 8
$ws.Range("AD8").Value = 6.8
$ws.Range("AE8").Value = 12.5
$ws.Range("AF8").Value = 50
$ws.Range("AG8").Value = 11.25
$ws.Range("AH8").Value = 18.5
$ws.Range("AI8").Value = 11.25
$ws.Range("AK8").Value = 27
$ws.Range("AL8").Value = 32
$ws.Range("AM8").Value = 350
$ws.Range("AN8").Value = 4.05
$ws.Range("AO8").Value = 10.25
$ws.Range("AP8").Value = 17
$ws.Range("AQ8").Value = 37
$ws.Range("AR8").Value = 60
$ws.Range("AS8").Value = 200
$ws.Range("AT8").Value = 2.9
$ws.Range("AU8").Value = 6.8
$ws.Range("AV8").Value = 55
$ws.Range("AX8").Value = 18
$ws.Range("AY8").Value = 23

# Row 9
$ws.Range("I9").Value = 2.5
$ws.Range("K9").Value = 2.1
$ws.Range("L9").Value = 3.1
$ws.Range("Q9").Value = 1.88
$ws.Range("T9").Value = 2.75
$ws.Range("AA9").Value = 21
$ws.Range("AB9").Value = 28
$ws.Range("AE9").Value = 13
$ws.Range("AH9").Value = 12.5
$ws.Range("AJ9").Value = 27
$ws.Range("AK9").Value = 20
$ws.Range("AT9").Value = 2.75
$ws.Range("AX9").Value = 13.5
$ws.Range("AY9").Value = 21
$ws.Range("BA9").Value = 90
